$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 736.9474
$ws.Range("I19").Value = 679.8570999999999
$ws.Range("J19").Value = 770.25
$ws.Range("K19").Value = 679.8570999999999
$ws.Range("L19").Value = 770.25
$ws.Range("M19").Value = -504.8570999999999
$ws.Range("N19").Value = -1120.25

# Row 38
$ws.Range("H38").Value = 578.4167
$ws.Range("I38").Value = 385.85715
$ws.Range("J38").Value = 848
$ws.Range("K38").Value = 1157.57145
$ws.Range("L38").Value = 2544
$ws.Range("M38").Value = -785.5714499999999
$ws.Range("N38").Value = -3288

# Row 39
$ws.Range("H39").Value = 155.31818
$ws.Range("I39").Value = 64.09090999999999
$ws.Range("J39").Value = 246.54546
$ws.Range("K39").Value = 192.27273
$ws.Range("L39").Value = 739.6363799999999
$ws.Range("M39").Value = 103.72727
$ws.Range("N39").Value = -1331.63638

# Row 40
$ws.Range("H40").Value = 2828.9285
$ws.Range("I40").Value = 3287.5
$ws.Range("J40").Value = 2217.5
$ws.Range("K40").Value = 3287.5
$ws.Range("L40").Value = 2217.5
$ws.Range("M40").Value = -3112.5
$ws.Range("N40").Value = -2567.5

# Row 42
$ws.Range("H42").Value = 82.55556
$ws.Range("I42").Value = 82.55556
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 247.66668
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -17.66667999999999

# Row 43
$ws.Range("H43").Value = 1152
$ws.Range("I43").Value = 1260
$ws.Range("J43").Value = 990
$ws.Range("K43").Value = 1260
$ws.Range("L43").Value = 990
$ws.Range("M43").Value = -1191
$ws.Range("N43").Value = -1128

# Row 97
$ws.Range("H97").Value = 800
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2400
$ws.Range("N97").Value = -3392

# Row 98
$ws.Range("H98").Value = 466979.47
$ws.Range("I98").Value = 559320.9
$ws.Range("J98").Value = 5272.5
$ws.Range("K98").Value = 559320.9
$ws.Range("L98").Value = 5272.5
$ws.Range("M98").Value = -557822.9
$ws.Range("N98").Value = -8268.5

# Row 101
$ws.Range("H101").Value = 13835.375
$ws.Range("I101").Value = 749
$ws.Range("J101").Value = 18197.5
$ws.Range("K101").Value = 2247
$ws.Range("L101").Value = 54592.5
$ws.Range("M101").Value = -625
$ws.Range("N101").Value = -57836.5

# Row 122
$ws.Range("H122").Value = 466979.47
$ws.Range("I122").Value = 559320.9
$ws.Range("J122").Value = 5272.5
$ws.Range("K122").Value = 1677962.7
$ws.Range("L122").Value = 15817.5
$ws.Range("M122").Value = -1675512.7
$ws.Range("N122").Value = -20717.5

# Row 133
$ws.Range("H133").Value = 18131.666
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 18131.666
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 18131.666
$ws.Range("N133").Value = -28251.666

# Row 138
$ws.Range("H138").Value = 8477135
$ws.Range("I138").Value = 2798.5
$ws.Range("J138").Value = 12822948
$ws.Range("K138").Value = 8395.5
$ws.Range("L138").Value = 38468844
$ws.Range("M138").Value = -3255.5
$ws.Range("N138").Value = -38479124

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2572.7212
$ws.Range("I32").Value = 2198.6155
$ws.Range("J32").Value = 3235.9092
$ws.Range("K32").Value = 2198.6155
$ws.Range("L32").Value = 3235.9092
$ws.Range("M32").Value = -1911.6155
$ws.Range("N32").Value = -3809.9092

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 133
$ws.Range("H133").Value = 42831.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 42831.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 42831.5
$ws.Range("N133").Value = -47891.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 15783.909
$ws.Range("I86").Value = 7361.6
$ws.Range("J86").Value = 100007
$ws.Range("K86").Value = 7361.6
$ws.Range("L86").Value = 100007
$ws.Range("M86").Value = -6238.6
$ws.Range("N86").Value = -102253

# Row 89
$ws.Range("H89").Value = 15783.909
$ws.Range("I89").Value = 7361.6
$ws.Range("J89").Value = 100007
$ws.Range("K89").Value = 36808
$ws.Range("L89").Value = 500035
$ws.Range("M89").Value = -31192
$ws.Range("N89").Value = -511267

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 1961.125
$ws.Range("I122").Value = 1670.6666
$ws.Range("J122").Value = 2135.4
$ws.Range("K122").Value = 5011.9998
$ws.Range("L122").Value = 6406.200000000001
$ws.Range("M122").Value = -2561.9998
$ws.Range("N122").Value = -11306.2

# Row 132
$ws.Range("H132").Value = 3131.7307
$ws.Range("I132").Value = 1760.2667
$ws.Range("J132").Value = 5001.909
$ws.Range("K132").Value = 5280.800099999999
$ws.Range("L132").Value = 15005.727
$ws.Range("M132").Value = -2750.800099999999
$ws.Range("N132").Value = -20065.727

# Row 134
$ws.Range("H134").Value = 6745.2856
$ws.Range("I134").Value = 4425.75
$ws.Range("J134").Value = 9838
$ws.Range("K134").Value = 13277.25
$ws.Range("L134").Value = 29514
$ws.Range("M134").Value = -10742.25
$ws.Range("N134").Value = -34584

$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 900
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2700
$ws.Range("N92").Value = -5196

# Row 104
$ws.Range("H104").Value = 10000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 10000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -35242
$ws.Range("M104").ClearContents()

# Row 123
$ws.Range("H123").Value = 2000
$ws.Range("I123").Value = 500
$ws.Range("J123").Value = 3500
$ws.Range("K123").Value = 1500
$ws.Range("L123").Value = 10500
$ws.Range("M123").Value = 950
$ws.Range("N123").Value = -15400

# Row 125
$ws.Range("H125").Value = 2999.4119
$ws.Range("I125").Value = 2995
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 8985
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -4065
$ws.Range("N125").Value = -18840

# Row 131
$ws.Range("H131").Value = 2804.1194
$ws.Range("I131").Value = 438.0909
$ws.Range("J131").Value = 3268.875
$ws.Range("K131").Value = 1314.2727
$ws.Range("L131").Value = 9806.625
$ws.Range("M131").Value = 3725.7273
$ws.Range("N131").Value = -19886.625

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 734.2632
$ws.Range("I22").Value = 770.9167
$ws.Range("J22").Value = 671.4286
$ws.Range("K22").Value = 770.9167
$ws.Range("L22").Value = 671.4286
$ws.Range("M22").Value = -475.9167
$ws.Range("N22").Value = -1261.4286

# Row 27
$ws.Range("H27").Value = 734.2632
$ws.Range("I27").Value = 770.9167
$ws.Range("J27").Value = 671.4286
$ws.Range("K27").Value = 770.9167
$ws.Range("L27").Value = 671.4286
$ws.Range("M27").Value = -663.9167
$ws.Range("N27").Value = -885.4286

$ws = $wb.Worksheets.Item("WVR")
# Row 114
$ws.Range("H114").Value = 100000
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 100000
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 100000
$ws.Range("N114").Value = -108678
